# Applies numeric corrections to Leve profit calculation sheets
# (currentAveragePrice / LevePrice / LeveProfit columns) across all
# job worksheets, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 395.48
$ws.Range("I15").Value = 395.48
$ws.Range("K15").Value = 1186.44
$ws.Range("M15").Value = -1017.44
$ws.Range("H40").Value = 7000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 7000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 7000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -7350
$ws.Range("H62").Value = 8122.0557
$ws.Range("I62").Value = 7042.7
$ws.Range("J62").Value = 9471.25
$ws.Range("K62").Value = 7042.7
$ws.Range("L62").Value = 9471.25
$ws.Range("M62").Value = -6418.7
$ws.Range("N62").Value = -10719.25
$ws.Range("H65").Value = 8122.0557
$ws.Range("I65").Value = 7042.7
$ws.Range("J65").Value = 9471.25
$ws.Range("K65").Value = 35213.5
$ws.Range("L65").Value = 47356.25
$ws.Range("M65").Value = -32093.5
$ws.Range("N65").Value = -53596.25
$ws.Range("H74").Value = 3750
$ws.Range("I74").Value = 3750
$ws.Range("K74").Value = 3750
$ws.Range("M74").Value = -2814
$ws.Range("H77").Value = 3750
$ws.Range("I77").Value = 3750
$ws.Range("K77").Value = 18750
$ws.Range("M77").Value = -14070
$ws.Range("H82").Value = 2104.25
$ws.Range("I82").Value = 2104.25
$ws.Range("K82").Value = 6312.75
$ws.Range("M82").Value = -5906.75
$ws.Range("H85").Value = 2104.25
$ws.Range("I85").Value = 2104.25
$ws.Range("K85").Value = 6312.75
$ws.Range("M85").Value = -4908.75
$ws.Range("H138").Value = 9381.35
$ws.Range("J138").Value = 9798.794
$ws.Range("L138").Value = 29396.382
$ws.Range("N138").Value = -39676.382
$ws.Range("H140").Value = 129998.2
$ws.Range("J140").Value = 129998.2
$ws.Range("L140").Value = 129998.2
$ws.Range("N140").Value = -140358.2

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2566.037
$ws.Range("I45").Value = 1168.3334
$ws.Range("K45").Value = 1168.3334
$ws.Range("M45").Value = -791.3334
$ws.Range("H132").Value = 3990.8696
$ws.Range("I132").Value = 3215.3845
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 9646.1535
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -7116.1535
$ws.Range("N132").Value = -20057
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 489.31818
$ws.Range("I22").Value = 489.78946
$ws.Range("J22").Value = 486.33334
$ws.Range("K22").Value = 489.78946
$ws.Range("L22").Value = 486.33334
$ws.Range("M22").Value = -316.78946
$ws.Range("N22").Value = -832.33334
$ws.Range("H132").Value = 191113.33
$ws.Range("J132").Value = 191113.33
$ws.Range("L132").Value = 191113.33
$ws.Range("N132").Value = -201233.33

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1000
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1224
$ws.Range("H31").Value = 5469.381
$ws.Range("I31").Value = 3799.5
$ws.Range("K31").Value = 3799.5
$ws.Range("M31").Value = -3504.5
$ws.Range("H34").Value = 5469.381
$ws.Range("I34").Value = 3799.5
$ws.Range("K34").Value = 3799.5
$ws.Range("M34").Value = -3597.5
$ws.Range("H122").Value = 2367.6
$ws.Range("I122").Value = 2297.3333
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6891.999899999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4441.999899999999
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 10988.167
$ws.Range("I132").Value = 12810.704
$ws.Range("J132").Value = 2786.75
$ws.Range("K132").Value = 38432.112
$ws.Range("L132").Value = 8360.25
$ws.Range("M132").Value = -35902.112
$ws.Range("N132").Value = -13420.25
$ws.Range("H134").Value = 3291.3142
$ws.Range("I134").Value = 2996.4482
$ws.Range("K134").Value = 8989.3446
$ws.Range("M134").Value = -6454.3446

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 96.666664
$ws.Range("I97").Value = 96.666664
$ws.Range("K97").Value = 289.999992
$ws.Range("M97").Value = 206.000008
$ws.Range("H122").Value = 2479.6155
$ws.Range("J122").Value = 2564.3333
$ws.Range("L122").Value = 23078.9997
$ws.Range("N122").Value = -27978.9997
$ws.Range("H132").Value = 2230.9443
$ws.Range("J132").Value = 2362.0625
$ws.Range("L132").Value = 21258.5625
$ws.Range("N132").Value = -26318.5625

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 12447
$ws.Range("J43").Value = 16155.333
$ws.Range("L43").Value = 16155.333
$ws.Range("N43").Value = -16457.333
$ws.Range("H46").Value = 54499.582
$ws.Range("J46").Value = 54499.582
$ws.Range("L46").Value = 54499.582
$ws.Range("N46").Value = -54811.582
$ws.Range("H55").Value = 11329.5
$ws.Range("I55").Value = 16000
$ws.Range("J55").Value = 10395.4
$ws.Range("K55").Value = 16000
$ws.Range("L55").Value = 10395.4
$ws.Range("M55").Value = -15673
$ws.Range("N55").Value = -11049.4
$ws.Range("H80").Value = 8597.200000000001
$ws.Range("J80").Value = 13529.4
$ws.Range("L80").Value = 13529.4
$ws.Range("N80").Value = -15525.4
$ws.Range("H83").Value = 8597.200000000001
$ws.Range("J83").Value = 13529.4
$ws.Range("L83").Value = 67647
$ws.Range("N83").Value = -77631

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 334093.34
$ws.Range("I2").Value = 730.43475
$ws.Range("J2").Value = 1429428.6
$ws.Range("K2").Value = 730.43475
$ws.Range("L2").Value = 1429428.6
$ws.Range("M2").Value = -618.43475
$ws.Range("N2").Value = -1429652.6
$ws.Range("H132").Value = 5805.4375
$ws.Range("I132").Value = 6474.3335
$ws.Range("K132").Value = 19423.0005
$ws.Range("M132").Value = -16893.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 12666.667
$ws.Range("J5").Value = 12666.667
$ws.Range("L5").Value = 12666.667
$ws.Range("N5").Value = -12890.667
$ws.Range("H132").Value = 6254738.5
$ws.Range("I132").Value = 7580713
$ws.Range("J132").Value = 3714.7144
$ws.Range("K132").Value = 22742139
$ws.Range("L132").Value = 11144.1432
$ws.Range("M132").Value = -22739609
$ws.Range("N132").Value = -16204.1432
$ws.Range("H136").Value = 13943.875
$ws.Range("I136").Value = 14337.662
$ws.Range("K136").Value = 43012.986
$ws.Range("M136").Value = -40462.986

